$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "#account.users"
$ws.Name = "#account.users"

# Widen column B to fit the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 27

# Shift header labels: drop the old "username" text from column B and
# replace it with the new clarified label; shift Ho/Ten/Mat khau left
# is not needed (they stay in C/D/E), but their shared-string text is
# re-set explicitly to make the resulting values match the target.
$ws.Range("B1").Value = "Tên đăng nhập (*, định danh)"
$ws.Range("C1").Value = "Họ (*)"
$ws.Range("D1").Value = "Tên (*)"
$ws.Range("E1").Value = "Mật khẩu (*)"

# Update the active selection to reflect the author's last cursor position
$ws.Range("B6").Select()
